$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected; unprotect to allow edits, then restore protection at the end.
$ws.Unprotect()

# Insert two blank rows before the old row 6 so the reference/header block
# (old rows 6-22) shifts down to rows 8-24, leaving a new gap at row 7
# (mirroring the original gap pattern used between sections).
$ws.Rows("6:7").Insert()

# Update the title block text.
$ws.Range("A2").Value = "Version 1.2.2"
$ws.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet."

# New instruction lines in the rows that were just inserted.
$ws.Range("A5").Value = "Do not change the headers of the 'Antibodies' sheet."
$ws.Range("A6").Value = "Do not edit the other sheets."

# Restore sheet protection to its original state.
$ws.Protect()
